# Scheduled-runner refresh of market-price-derived leve profit columns
# (currentAveragePrice / NQ / HQ / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 398
$ws.Range("I19").Value = 411.0909
$ws.Range("J19").Value = 382
$ws.Range("K19").Value = 411.0909
$ws.Range("L19").Value = 382
$ws.Range("M19").Value = -236.0909
$ws.Range("N19").Value = -732

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 466.53333
$ws.Range("I80").Value = 261.85715
$ws.Range("J80").Value = 645.625
$ws.Range("K80").Value = 785.5714499999999
$ws.Range("L80").Value = 1936.875
$ws.Range("M80").Value = 212.4285500000001
$ws.Range("N80").Value = -3932.875

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 466.53333
$ws.Range("I83").Value = 261.85715
$ws.Range("J83").Value = 645.625
$ws.Range("K83").Value = 2356.71435
$ws.Range("L83").Value = 5810.625
$ws.Range("M83").Value = 2635.28565
$ws.Range("N83").Value = -15794.625

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 3035.25
$ws.Range("I111").Value = 2984.8333
$ws.Range("K111").Value = 8954.499899999999
$ws.Range("M111").Value = -5887.499899999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 3405.375
$ws.Range("I131").Value = 3463.2856
$ws.Range("J131").Value = 3000
$ws.Range("K131").Value = 10389.8568
$ws.Range("L131").Value = 9000
$ws.Range("M131").Value = -5349.856800000001
$ws.Range("N131").Value = -19080

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2075.1396
$ws.Range("I132").Value = 1964.0264
$ws.Range("J132").Value = 2919.6
$ws.Range("K132").Value = 5892.0792
$ws.Range("L132").Value = 8758.799999999999
$ws.Range("M132").Value = -3362.0792
$ws.Range("N132").Value = -13818.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 7896.3
$ws.Range("I137").Value = 13393.6
$ws.Range("K137").Value = 40180.8
$ws.Range("M137").Value = -37630.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 10384.6
$ws.Range("I141").Value = 7980.75
$ws.Range("K141").Value = 23942.25
$ws.Range("M141").Value = -18762.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4285.387
$ws.Range("I32").Value = 3826
$ws.Range("K32").Value = 3826
$ws.Range("M32").Value = -3539

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 40949.5
$ws.Range("I34").Value = 40949.5
$ws.Range("K34").Value = 40949.5
$ws.Range("M34").Value = -40678.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2126.0488
$ws.Range("I61").Value = 1477.0834
$ws.Range("J61").Value = 6798.6
$ws.Range("K61").Value = 1477.0834
$ws.Range("L61").Value = 6798.6
$ws.Range("M61").Value = -1265.0834
$ws.Range("N61").Value = -7222.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3561.303
$ws.Range("I74").Value = 3189.4348
$ws.Range("K74").Value = 3189.4348
$ws.Range("M74").Value = -2315.4348

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3561.303
$ws.Range("I77").Value = 3189.4348
$ws.Range("K77").Value = 15947.174
$ws.Range("M77").Value = -11579.174

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3990.6086
$ws.Range("I122").Value = 4123.421
$ws.Range("J122").Value = 3359.75
$ws.Range("K122").Value = 12370.263
$ws.Range("L122").Value = 10079.25
$ws.Range("M122").Value = -9920.263000000001
$ws.Range("N122").Value = -14979.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2126.0488
$ws.Range("I136").Value = 1477.0834
$ws.Range("J136").Value = 6798.6
$ws.Range("K136").Value = 4431.2502
$ws.Range("L136").Value = 20395.8
$ws.Range("M136").Value = -1881.2502
$ws.Range("N136").Value = -25495.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3479.762
$ws.Range("I134").Value = 2218.8276
$ws.Range("K134").Value = 6656.4828
$ws.Range("M134").Value = -4121.4828

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2174.4
$ws.Range("I31").Value = 1651.6
$ws.Range("K31").Value = 1651.6
$ws.Range("M31").Value = -1356.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2174.4
$ws.Range("I34").Value = 1651.6
$ws.Range("K34").Value = 1651.6
$ws.Range("M34").Value = -1449.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 43525.75
$ws.Range("J59").Value = 56999.5
$ws.Range("L59").Value = 56999.5
$ws.Range("N59").Value = -59289.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1362.3448
$ws.Range("I122").Value = 1104.1364
$ws.Range("K122").Value = 3312.4092
$ws.Range("M122").Value = -862.4092000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1449.8
$ws.Range("I132").Value = 1235.6923
$ws.Range("K132").Value = 3707.0769
$ws.Range("M132").Value = -1177.0769

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1577.1111
$ws.Range("I134").Value = 1564.1613
$ws.Range("K134").Value = 4692.4839
$ws.Range("M134").Value = -2157.4839

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 94433.414
$ws.Range("I4").Value = 125449.875
$ws.Range("K4").Value = 376349.625
$ws.Range("M4").Value = -376237.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1405.4736
$ws.Range("I5").Value = 1097.4286
$ws.Range("J5").Value = 2268
$ws.Range("K5").Value = 3292.2858
$ws.Range("L5").Value = 6804
$ws.Range("M5").Value = -3180.2858
$ws.Range("N5").Value = -7028

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 53.666668
$ws.Range("I12").Value = 52
$ws.Range("J12").Value = 56
$ws.Range("K12").Value = 156
$ws.Range("L12").Value = 168
$ws.Range("M12").Value = 17
$ws.Range("N12").Value = -514

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1039
$ws.Range("J98").Value = 1048.75
$ws.Range("L98").Value = 3146.25
$ws.Range("N98").Value = -6142.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 443.55554
$ws.Range("I107").Value = 434.6
$ws.Range("K107").Value = 1303.8
$ws.Range("M107").Value = 616.1999999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 59169.223
$ws.Range("I117").Value = 435.25
$ws.Range("J117").Value = 75950.36
$ws.Range("K117").Value = 1305.75
$ws.Range("L117").Value = 227851.08
$ws.Range("M117").Value = 2136.25
$ws.Range("N117").Value = -234735.08

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2423.75
$ws.Range("J132").Value = 2498.3333
$ws.Range("L132").Value = 22484.9997
$ws.Range("N132").Value = -27544.9997

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1405.4736
$ws.Range("I135").Value = 1097.4286
$ws.Range("J135").Value = 2268
$ws.Range("K135").Value = 9876.857399999999
$ws.Range("L135").Value = 20412
$ws.Range("M135").Value = -7341.857399999999
$ws.Range("N135").Value = -25482

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 37507.25
$ws.Range("J49").Value = 37507.25
$ws.Range("L49").Value = 37507.25
$ws.Range("N49").Value = -37875.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3799.875
$ws.Range("I80").Value = 3235.6667
$ws.Range("J80").Value = 4138.4
$ws.Range("K80").Value = 3235.6667
$ws.Range("L80").Value = 4138.4
$ws.Range("M80").Value = -2237.6667
$ws.Range("N80").Value = -6134.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3799.875
$ws.Range("I83").Value = 3235.6667
$ws.Range("J83").Value = 4138.4
$ws.Range("K83").Value = 16178.3335
$ws.Range("L83").Value = 20692
$ws.Range("M83").Value = -11186.3335
$ws.Range("N83").Value = -30676

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 57466
$ws.Range("I42").Value = 57466
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 57466
$ws.Range("L42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("N42").Value = -56903

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H49").Value = 57466
$ws.Range("I49").Value = 57466
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 57466
$ws.Range("L49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -57319

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3323.5
$ws.Range("I68").Value = 1188.2
$ws.Range("J68").Value = 14000
$ws.Range("K68").Value = 1188.2
$ws.Range("L68").Value = 14000
$ws.Range("M68").Value = -439.2
$ws.Range("N68").Value = -15498

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 3323.5
$ws.Range("I71").Value = 1188.2
$ws.Range("J71").Value = 14000
$ws.Range("K71").Value = 5941
$ws.Range("L71").Value = 70000
$ws.Range("M71").Value = -2197
$ws.Range("N71").Value = -77488

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2293.7896
$ws.Range("I82").Value = 1083.75
$ws.Range("J82").Value = 3173.818
$ws.Range("K82").Value = 1083.75
$ws.Range("L82").Value = 3173.818
$ws.Range("M82").Value = -722.75
$ws.Range("N82").Value = -3895.818

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2293.7896
$ws.Range("I85").Value = 1083.75
$ws.Range("J85").Value = 3173.818
$ws.Range("K85").Value = 1083.75
$ws.Range("L85").Value = 3173.818
$ws.Range("M85").Value = 164.25
$ws.Range("N85").Value = -5669.818

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3318.6882
$ws.Range("I132").Value = 2884.6584
$ws.Range("J132").Value = 6554.1816
$ws.Range("K132").Value = 8653.975199999999
$ws.Range("L132").Value = 19662.5448
$ws.Range("M132").Value = -6123.975199999999
$ws.Range("N132").Value = -24722.5448

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 28789.533
$ws.Range("I62").Value = 18570.584
$ws.Range("J62").Value = 69665.336
$ws.Range("K62").Value = 18570.584
$ws.Range("L62").Value = 69665.336
$ws.Range("M62").Value = -17946.584
$ws.Range("N62").Value = -70913.336

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 28789.533
$ws.Range("I65").Value = 18570.584
$ws.Range("J65").Value = 69665.336
$ws.Range("K65").Value = 92852.92
$ws.Range("L65").Value = 348326.68
$ws.Range("M65").Value = -89732.92
$ws.Range("N65").Value = -354566.68

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 927.08105
$ws.Range("J132").Value = 2495
$ws.Range("L132").Value = 7485
$ws.Range("N132").Value = -12545
